# Fixed the insert-rows and column maping
#
# The Phone column (D) values were missing the "801-" area code prefix.
# Update each associate's phone number in place (D2:D7) so the shared
# strings table is rewritten the same way Excel does it: the old
# "555-XXXX" strings drop out of use and new "801-555-XXXX" strings get
# appended. Associates are touched in (Jane, John, Kenji, Maria) order -
# i.e. row 7 (Kenji) is written before row 6 (Maria) - to match how the
# workbook's shared-string table ends up ordered after the real edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "801-555-1234"
$ws.Range("D3").Value = "801-555-1234"
$ws.Range("D4").Value = "801-555-1234"
$ws.Range("D5").Value = "801-555-5678"
$ws.Range("D7").Value = "801-555-3456"
$ws.Range("D6").Value = "801-555-9012"

# Restore the view state: active cell / selection moved to J7, and the
# window zoomed to 85% so the job columns (H:M) are in frame.
$excel.ActiveWindow.Zoom = 85
$ws.Range("J7").Select()
